# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns to the
# latest scrape. D-column values are forced to Text (via a leading
# apostrophe, matching how Excel's own "force text" entry works) so that
# digits like trailing zeros / dotted thousand-separators ("0.420",
# "66.513.09") survive instead of being auto-coerced into numbers; the
# Style is then reset to "Normal" so no stray text-format style sticks to
# the cell (keeping it identical to the untouched cells around it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.513.09"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.83%  '
$ws.Range('D3').Value = "'3.433.80"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'580.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.09%  '
$ws.Range('D6').Value = "'174.51"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.46%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'0.599"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('D9').Value = "'3.430.38"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('D10').Value = "'0.134"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.93%  '
$ws.Range('D11').Value = "'6.89"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.38%  '
$ws.Range('D12').Value = "'0.420"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.36%  '
$ws.Range('D13').Value = "'4.033.01"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('D14').Value = "'30.85"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.76%  '
$ws.Range('D15').Value = "'0.131"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.73%  '
$ws.Range('D16').Value = "'66.527.01"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('D17').Value = "'0.0000172"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.08%  '
$ws.Range('D18').Value = "'3.445.00"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('D19').Value = "'6.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.90%  '
$ws.Range('D20').Value = "'13.73"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.82%  '
$ws.Range('D21').Value = "'375.53"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.42%  '
$ws.Range('D22').Value = "'7.75"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.31%  '
$ws.Range('D23').Value = "'0.998"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = "'5.72"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = "'71.01"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.26%  '
$ws.Range('D26').Value = "'0.526"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('D27').Value = "'0.0000118"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.28%  '
$ws.Range('D28').Value = "'9.80"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.97%  '
$ws.Range('D29').Value = "'0.172"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.88%  '
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('D31').Value = "'5.87"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.57%  '
$ws.Range('D32').Value = "'2.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.12%  '
$ws.Range('D33').Value = "'23.75"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  -6.81%  '
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('D36').Value = "'7.11"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.41%  '
$ws.Range('D37').Value = "'1.52"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.06%  '
$ws.Range('D38').Value = "'160.36"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.34%  '
$ws.Range('D39').Value = "'0.873"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').Value = "'26.93"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.68%  '
$ws.Range('D41').Value = "'1.79"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.29%  '
$ws.Range('D42').Value = "'2.60"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.35%  '
$ws.Range('D43').Value = "'6.54"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.92%  '
$ws.Range('D44').Value = "'4.43"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.53%  '
$ws.Range('D45').Value = "'2.682.69"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.33%  '
$ws.Range('D46').Value = "'0.0689"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.73%  '
$ws.Range('D47').Value = "'25.21"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.26%  '
$ws.Range('D48').Value = "'40.79"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.40%  '
$ws.Range('D49').Value = "'0.0292"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.46%  '
$ws.Range('D50').Value = "'315.08"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.97%  '
$ws.Range('D51').Value = "'1.01"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.30%  '
